$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns (L:S) -------------------------------------------
$headers = @(
    "Age",
    "Gender",
    "Nationality",
    "Passport Number",
    "Passport Validity",
    "NOK Name",
    "NOK Contact Number",
    "Medical Clearances"
)

$col = 12
foreach ($h in $headers) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $h
    $cell.WrapText = $true
    $col = $col + 1
}

# --- Row 1 height (75 -> 45) ---------------------------------------------
$ws.Rows.Item(1).RowHeight = 45

# --- Column widths ----------------------------------------------------
# (ColumnWidth values chosen so the saved XML width lands as close as
# possible to the widths recorded by the source workbook.)
$ws.Columns.Item(1).ColumnWidth  = 12.6666666666667
$ws.Columns.Item(2).ColumnWidth  = 12.3333333333333
$ws.Columns.Item(3).ColumnWidth  = 11.5
$ws.Columns.Item(4).ColumnWidth  = 14.6666666666667
$ws.Columns.Item(5).ColumnWidth  = 11.3333333333333
$ws.Columns.Item(7).ColumnWidth  = 10.1666666666667
$ws.Columns.Item(8).ColumnWidth  = 10.8333333333333
$ws.Columns.Item(10).ColumnWidth = 9.83333333333333
$ws.Columns.Item(12).ColumnWidth = 8.83333333333333
$ws.Columns.Item(13).ColumnWidth = 11.6666666666667
$ws.Columns.Item(14).ColumnWidth = 14
$ws.Columns.Item(15).ColumnWidth = 15.3333333333333
$ws.Columns.Item(16).ColumnWidth = 14.6666666666667
$ws.Columns.Item(17).ColumnWidth = 17.6666666666667
$ws.Columns.Item(18).ColumnWidth = 13.6666666666667
$ws.Columns.Item(19).ColumnWidth = 17

# --- View / selection ----------------------------------------------------
$ws.Range("S1").Select() | Out-Null
